# Added +1 hours for me & Matt
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Andrew Case: 59h 35m -> 60h 35m (+1 hour)
$ws.Range("B4").Value = "60h 35m"

# Matthew Darby: 28h 50m -> 29h 50m (+1 hour)
$ws.Range("B5").Value = "29h 50m"

# Move the active selection to B5, matching the author's edit location
$ws.Range("B5").Select()
